$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new invoice rows (3-5) ---
$ws.Range("A3").Value = "INV-10012"
$ws.Range("B3").Value = 44281
$ws.Range("C3").Value = 44311
$ws.Range("D3").Value = 1699.48

$ws.Range("A4").Value = 202205
$ws.Range("B4").Value = 44712
$ws.Range("C4").Value = 44726
$ws.Range("D4").Value = 220

$ws.Range("A5").Value = "A246"
$ws.Range("B5").Value = 43174
$ws.Range("C5").Value = 43205
$ws.Range("D5").Value = 700.65

# Reuse the existing date number formatting (style applied to B2:C2) for the
# newly added date cells, so no redundant style/numFmt gets created.
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B3:C5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Swap theme accent1 / accent5 colors ---
$theme = $wb.Theme
$accent1 = $theme.ThemeColorScheme(5)
$accent5 = $theme.ThemeColorScheme(9)
$accent1Rgb = $accent1.RGB()
$accent5Rgb = $accent5.RGB()
$accent1.RGB = $accent5Rgb
$accent5.RGB = $accent1Rgb
